$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auto-updated maa:// operator-recognition accuracy figures (CI: Auto Update Data)
$ws.Range("D2").Value = "maa://24702 (94.29), maa://25390 (95.98), maa://36681 (87.01)"
$ws.Range("L2").Value = "*maa://24633 (56.05), *maa://30515 (69.31), *maa://34787 (72.86), ***maa://20792 (11.93), maa://39402 (89.13), ***maa://29083 (27.78)"
$ws.Range("X3").Value = "maa://27396 (84.57), maa://27484 (96.15), maa://27480 (82.86)"
$ws.Range("X4").Value = "**maa://32495 (47.91), ***maa://31785 (22.22), ***maa://36683 (28.26), maa://43217 (90.62)"
$ws.Range("D5").Value = "maa://21245 (83.56), maa://22744 (84.0)"
$ws.Range("AF6").Value = "*maa://33152 (60.0), ***maa://22770 (26.09)"
$ws.Range("L7").Value = "maa://28624 (93.55), maa://24957 (97.62)"
$ws.Range("X7").Value = "maa://22399 (95.17), *maa://22758 (74.19)"
$ws.Range("A8").Value = "更新日期：2025.01.04 13:16:59"
$ws.Range("D8").Value = "*maa://21476 (72.92), **maa://39431 (50.0), *maa://37551 (57.14)"
$ws.Range("AF8").Value = "*maa://24479 (76.83), *maa://21990 (51.85)"
$ws.Range("L9").Value = "maa://22762 (91.95), maa://39552 (90.0)"
$ws.Range("AF9").Value = "maa://26206 (89.42), *maa://22865 (51.92)"
$ws.Range("T10").Value = "maa://27395 (96.11), maa://22755 (87.61), **maa://22756 (40.91), ***maa://21737 (10.61)"
$ws.Range("AE11").Value = "1"
$ws.Range("AF11").Value = "maa://31203 (95.65)"
$ws.Range("X12").Value = "maa://22753 (91.46), *maa://21485 (77.04), maa://37962 (89.66)"
$ws.Range("D13").Value = "maa://24999 (91.8), maa://36673 (92.65), maa://25001 (85.51)"
$ws.Range("X13").Value = "*maa://34957 (80.0), *maa://22768 (51.61)"
$ws.Range("L14").Value = "maa://26245 (96.55), maa://21288 (96.3), maa://39841 (95.29), maa://36682 (97.37)"
$ws.Range("D15").Value = "*maa://22743 (77.44), maa://22734 (84.03), *maa://30808 (65.08), **maa://36048 (34.04), maa://45058 (100.0)"
$ws.Range("T16").Value = "maa://22729 (94.77), *maa://28648 (68.85), maa://36674 (83.33)"
$ws.Range("D18").Value = "maa://24570 (97.14)"
$ws.Range("AB19").Value = "*maa://30709 (63.59), *maa://36668 (56.41)"
$ws.Range("D20").Value = "maa://21432 (89.93), maa://25198 (93.14), *maa://20795 (51.18), maa://36680 (93.55)"
$ws.Range("L20").Value = "maa://41331 (84.62)"
$ws.Range("AB21").Value = "maa://21443 (80.23), ***maa://23820 (29.82)"
$ws.Range("AF21").Value = "maa://22524 (94.63), *maa://22432 (77.05)"
$ws.Range("L22").Value = "maa://27127 (85.58), *maa://22751 (73.85)"
$ws.Range("L23").Value = "maa://39756 (94.36), maa://39875 (93.85)"
$ws.Range("D24").Value = "*maa://24368 (78.3)"
$ws.Range("X24").Value = "maa://29988 (86.46), maa://23504 (93.18), **maa://22892 (39.31), *maa://25141 (76.38), *maa://36663 (78.26), ***maa://22815 (23.08)"
$ws.Range("H25").Value = "*maa://29063 (74.0), *maa://25311 (73.53), ***maa://22725 (4.84), maa://45047 (100.0)"
$ws.Range("AB25").Value = "maa://31215 (86.0), *maa://24516 (79.78), maa://26001 (87.5)"
$ws.Range("AB26").Value = "maa://42235 (93.51)"
$ws.Range("X28").Value = "maa://39929 (89.85), ***maa://39723 (14.29), maa://41749 (91.67)"
$ws.Range("AF28").Value = "maa://36660 (93.04), *maa://36701 (64.29)"
$ws.Range("H29").Value = "*maa://25175 (68.75)"
$ws.Range("L29").Value = "maa://28432 (92.95), *maa://28440 (76.84), maa://31400 (100.0), *maa://28650 (71.43)"
$ws.Range("AF29").Value = "*maa://24080 (69.05), *maa://42865 (78.95), ***maa://34960 (8.33)"
$ws.Range("AB30").Value = "maa://42979 (96.52), maa://45045 (100.0)"
$ws.Range("L31").Value = "maa://35926 (93.73), maa://36258 (83.67), *maa://43904 (80.0)"
$ws.Range("H32").Value = "maa://21895 (97.37), maa://36667 (98.41), **maa://20793 (38.78), maa://22760 (100.0)"
$ws.Range("T32").Value = "maa://42859 (96.47), maa://41108 (87.76), maa://41238 (96.3)"
$ws.Range("L35").Value = "maa://41296 (96.75)"
$ws.Range("AF38").Value = "maa://36697 (86.41)"
$ws.Range("H39").Value = "maa://25199 (84.82), maa://36670 (87.78), maa://30434 (89.39), ***maa://25036 (16.0), *maa://44165 (66.67), *maa://45059 (66.67)"
$ws.Range("T45").Value = "**maa://39364 (36.36)"
$ws.Range("H46").Value = "maa://35931 (92.63), maa://43901 (88.89)"
$ws.Range("H53").Value = "maa://32534 (93.6), **maa://32434 (34.78)"
$ws.Range("H55").Value = "maa://32532 (91.98)"
$ws.Range("H59").Value = "maa://27746 (83.18), maa://31270 (94.92)"
$ws.Range("H62").Value = "maa://42981 (96.77), maa://43903 (100.0)"
